# Refresh the cryptocurrency Price (column D) and Volume(1h) (column E)
# figures with the latest scraped values (GitHub Actions data refresh).
#
# Every cell in this sheet is stored as text. Price values that look like
# plain decimal numbers (e.g. "229.42") are written with a leading
# apostrophe so Excel keeps them as text instead of auto-converting them
# to numbers, matching the rest of the data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.120.01"
$ws.Range("E2").Value = "  +10.61%  "
$ws.Range("D3").Value = "1.819.51"
$ws.Range("E3").Value = "  +7.47%  "
$ws.Range("E4").Value = "  +0.53%  "
$ws.Range("D5").Value = "'229.42"
$ws.Range("E5").Value = "  +3.36%  "
$ws.Range("D6").Value = "'0.542"
$ws.Range("E6").Value = "  +3.77%  "
$ws.Range("E7").Value = "  +0.46%  "
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("D9").Value = "'45.92"
$ws.Range("E9").Value = "  +3.26%  "
$ws.Range("E10").Value = "  +5.47%  "
$ws.Range("E11").Value = "  +7.22%  "
$ws.Range("D12").Value = "'0.0930"
$ws.Range("D13").Value = "2.080.15"
$ws.Range("E13").Value = "  +7.42%  "
$ws.Range("D14").Value = "1.820.05"
$ws.Range("E14").Value = "  +7.35%  "
$ws.Range("E15").Value = "  +2.75%  "
$ws.Range("D16").Value = "34.117.35"
$ws.Range("E16").Value = "  +10.57%  "
$ws.Range("D17").Value = "'10.23"
$ws.Range("E17").Value = "  -4.29%  "
$ws.Range("D18").Value = "'4.32"
$ws.Range("E18").Value = "  +7.38%  "
$ws.Range("D19").Value = "'69.69"
$ws.Range("E19").Value = "  +4.72%  "
$ws.Range("D20").Value = "'258.04"
$ws.Range("E20").Value = "  +3.68%  "
$ws.Range("D21").Value = "0.0₃0749"
$ws.Range("E21").Value = "  +3.94%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.62%  "
$ws.Range("D23").Value = "'10.46"
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("D24").Value = "'4.36"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("D25").Value = "'2.19"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").Value = "'161.29"
$ws.Range("E26").Value = "  +2.43%  "
$ws.Range("D27").Value = "'16.73"
$ws.Range("E27").Value = "  +4.88%  "
$ws.Range("E28").Value = "  +5.69%  "
$ws.Range("E29").Value = "  +3.76%  "
$ws.Range("E30").Value = "  +0.52%  "
$ws.Range("D31").Value = "'3.89"
$ws.Range("E31").Value = "  +11.21%  "
$ws.Range("D32").Value = "'0.0515"
$ws.Range("E32").Value = "  +2.60%  "
$ws.Range("D33").Value = "'1.21"
$ws.Range("E33").Value = "  +5.85%  "
$ws.Range("D34").Value = "'3.57"
$ws.Range("E34").Value = "  +7.59%  "
$ws.Range("D35").Value = "1.569.13"
$ws.Range("E35").Value = "  +2.98%  "
$ws.Range("D36").Value = "'1.83"
$ws.Range("E36").Value = "  +4.16%  "
$ws.Range("E37").Value = "  +3.25%  "
$ws.Range("E38").Value = "  +4.84%  "
$ws.Range("D39").Value = "'84.76"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").Value = "'0.626"
$ws.Range("E40").Value = "  +6.62%  "
$ws.Range("D41").Value = "'2.84"
$ws.Range("E41").Value = "  +4.24%  "
$ws.Range("E42").Value = "  +0.62%  "
$ws.Range("D43").Value = "'0.911"
$ws.Range("E43").Value = "  +6.39%  "
$ws.Range("E44").Value = "  +5.00%  "
$ws.Range("D45").Value = "'0.0520"
$ws.Range("E45").Value = "  +3.24%  "
$ws.Range("E46").Value = "  +4.52%  "
$ws.Range("D47").Value = "1.974.48"
$ws.Range("E47").Value = "  +7.68%  "
$ws.Range("D48").Value = "'5.71"
$ws.Range("E48").Value = "  +4.87%  "
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("D50").Value = "'52.98"
$ws.Range("E50").Value = "  +0.90%  "
$ws.Range("E51").Value = "  +8.93%  "
